# MAJ plan afficheur temp
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update the "Afficheur Temp / Volt" segment/digit mapping table (columns I:J)
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 2

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 2

$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 7

$ws.Range("J18").Value = 3

$ws.Range("J19").Value = 4

# Update the active selection shown when the sheet is reopened
$ws.Range("N13").Select()
